# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46003
$ws.Range("B2").Value = 94.05
$ws.Range("C2").Value = 79.06
$ws.Range("D2").Value = 76.78
$ws.Range("E2").Value = 75.23999999999999
$ws.Range("F2").Value = 75.70999999999999
$ws.Range("G2").Value = 82.27
$ws.Range("H2").Value = 94.44
$ws.Range("I2").Value = 104.19
$ws.Range("J2").Value = 119.18
$ws.Range("K2").Value = 110.19
$ws.Range("L2").Value = 102.93
$ws.Range("M2").Value = 103.63
$ws.Range("N2").Value = 100.19
$ws.Range("O2").Value = 98.15000000000001
$ws.Range("P2").Value = 99.11
$ws.Range("Q2").Value = 94.44
$ws.Range("R2").Value = 97.27
$ws.Range("S2").Value = 105.05
$ws.Range("T2").Value = 112.98
$ws.Range("U2").Value = 109.45
$ws.Range("V2").Value = 108.99
$ws.Range("W2").Value = 104.82
$ws.Range("X2").Value = 96.13
$ws.Range("Y2").Value = 88.23999999999999
$ws.Range("Z2").Value = 97.19
$ws.Range("AA2").Value = "8h-12h"
$ws.Range("AB2").Value = 108.98
$ws.Range("AC2").Value = "8h-10h"
$ws.Range("AD2").Value = 114.68
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 111.22
$ws.Range("AG2").Value = "0h-23h"
